{"js": "// Remove the blank \"CV - Experience\" styled paragraph that sits between\n// the \"{years_exp}\" line and the \"Profile\" title (blank line between\n// profile and years of experience).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.style === \"CV - Experience\" && para.text.trim() === \"\") {\n    para.delete();\n    break;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the blank line (an empty \"CV - Experience\" styled paragraph)\n# that sits between the years-of-experience line and the \"Profile\"\n# title (\"Removed Blank line between profile and years of experience\").\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Style.NameLocal -eq \"CV - Experience\" -and $p.Range.Text.Trim() -eq \"\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $target.Range.Delete()\n}\n"}
